# In PowerPoint, the content of a top-level list is at the same level as
# the content of a top-level paragraph - the only difference is that a
# list style has been applied.
#
# The pptx writer used to increment the paragraph level on each list,
# turning what should be top-level lists into second-level lists. This
# script fixes that by walking every paragraph in every text frame (on
# every slide) and, wherever the paragraph sits at the (incorrectly
# bumped) second outline level, pulling it back up to the top level.
#
# COM's TextRange.IndentLevel is 1-based (IndentLevel = 1 corresponds to
# <a:pPr lvl="0"/> in the OOXML, IndentLevel = 2 corresponds to
# <a:pPr lvl="1"/>, and so on), so "second level" paragraphs have
# IndentLevel -eq 2.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $textRange = $shape.TextFrame.TextRange
        $paraCount = $textRange.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $textRange.Paragraphs($pi, 1)

            if ($para.IndentLevel -eq 2) {
                $para.IndentLevel = 1
            }
        }
    }
}
